$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so Excel does not
# auto-convert them to numbers (which would also lose exact formatting).
$textCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D16", "D17", "D19", "D21", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D43", "D45", "D46", "D47", "D48", "D50")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values
$ws.Range('D2').Value = '43.930.36'
$ws.Range('E2').Value = '  -5.58%  '
$ws.Range('D3').Value = '2.583.23'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '302.34'
$ws.Range('E5').Value = '  -1.95%  '
$ws.Range('D6').Value = '96.76'
$ws.Range('E6').Value = '  -4.07%  '
$ws.Range('D7').Value = '0.582'
$ws.Range('E7').Value = '  -3.28%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = '0.565'
$ws.Range('E9').Value = '  -2.09%  '
$ws.Range('D10').Value = '37.10'
$ws.Range('E10').Value = '  -6.06%  '
$ws.Range('D11').Value = '0.0820'
$ws.Range('E11').Value = '  -3.19%  '
$ws.Range('D12').Value = '7.85'
$ws.Range('E12').Value = '  -4.01%  '
$ws.Range('D13').Value = '2.974.55'
$ws.Range('E13').Value = '  -0.97%  '
$ws.Range('E14').Value = '  +1.33%  '
$ws.Range('D15').Value = '2.578.86'
$ws.Range('E15').Value = '  -1.03%  '
$ws.Range('D16').Value = '0.895'
$ws.Range('E16').Value = '  -2.79%  '
$ws.Range('D17').Value = '14.43'
$ws.Range('E17').Value = '  -3.68%  '
$ws.Range('D18').Value = '43.886.05'
$ws.Range('E18').Value = '  -5.80%  '
$ws.Range('D19').Value = '6.73'
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('D20').Value = '0.0₃0987'
$ws.Range('E20').Value = '  -2.79%  '
$ws.Range('D21').Value = '12.55'
$ws.Range('E21').Value = '  -3.40%  '
$ws.Range('E22').Value = '  +2.96%  '
$ws.Range('D23').Value = '266.84'
$ws.Range('E23').Value = '  -3.09%  '
$ws.Range('B24').Value = 'ImmutableX'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D24').Value = '2.24'
$ws.Range('E24').Value = '  +3.37%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = '2.95'
$ws.Range('E25').Value = '  -2.47%  '
$ws.Range('D26').Value = '29.51'
$ws.Range('E26').Value = '  +1.93%  '
$ws.Range('D28').Value = '10.29'
$ws.Range('E28').Value = '  -3.06%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').Value = '38.42'
$ws.Range('E29').Value = '  -1.87%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '2.16'
$ws.Range('E30').Value = '  -5.91%  '
$ws.Range('D31').Value = '6.24'
$ws.Range('E31').Value = '  -1.93%  '
$ws.Range('D32').Value = '3.60'
$ws.Range('E32').Value = '  -0.69%  '
$ws.Range('D33').Value = '2.22'
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').Value = '152.92'
$ws.Range('E34').Value = '  +1.05%  '
$ws.Range('B35').Value = 'WEMIXToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = '2.80'
$ws.Range('E35').Value = '  -2.22%  '
$ws.Range('D36').Value = '0.0817'
$ws.Range('E36').Value = '  -2.85%  '
$ws.Range('D37').Value = '0.118'
$ws.Range('E37').Value = '  -4.15%  '
$ws.Range('E38').Value = '  -1.60%  '
$ws.Range('D39').Value = '23.88'
$ws.Range('E39').Value = '  +2.22%  '
$ws.Range('D40').Value = '16.89'
$ws.Range('E40').Value = '  +5.98%  '
$ws.Range('D41').Value = '3.59'
$ws.Range('E41').Value = '  -1.81%  '
$ws.Range('E42').Value = '  -4.67%  '
$ws.Range('D43').Value = '3.90'
$ws.Range('E43').Value = '  -4.60%  '
$ws.Range('D44').Value = '2.042.43'
$ws.Range('E44').Value = '  -4.56%  '
$ws.Range('D45').Value = '0.997'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('D46').Value = '88.23'
$ws.Range('E46').Value = '  -5.22%  '
$ws.Range('D47').Value = '9.18'
$ws.Range('E47').Value = '  -3.30%  '
$ws.Range('D48').Value = '1.64'
$ws.Range('E48').Value = '  +6.28%  '
$ws.Range('D49').Value = '2.831.76'
$ws.Range('E49').Value = '  -0.82%  '
$ws.Range('D50').Value = '105.91'
$ws.Range('E50').Value = '  -2.73%  '
$ws.Range('E51').Value = '  -3.96%  '
